$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns: A=audit_id, B=name, C=type, D=major
# Target layout: A=name, B=type, C=major, D=audit_id
# i.e. each row's 4 values are rotated left by one column (A->D, B->A, C->B, D->C)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 1).Value2 = $b
    $ws.Cells.Item($r, 2).Value2 = $c
    $ws.Cells.Item($r, 3).Value2 = $d
    $ws.Cells.Item($r, 4).Value2 = $a
}
